{"js": "// Load all paragraphs in the document body so we can locate the ones we\n// need to edit by their current text.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet gtPara = null;      // \">>  >  your stuff after this line >>>\" paragraph\nlet bazPara = null;     // \"Baz chan\" + bookmark + \"ges\" paragraph\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"your stuff after this line\") !== -1) {\n    gtPara = paragraphs.items[i];\n  } else if (t.indexOf(\"Baz chan\") !== -1) {\n    bazPara = paragraphs.items[i];\n  }\n}\n\n// 1) Normalize the \">>>  your stuff after this line >>>\" paragraph into a\n//    single run (removes the stray proofing-error split without changing\n//    the visible text).\nif (gtPara) {\n  gtPara.insertText(\">>>  your stuff after this line >>>\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n\n// 2) Turn the \"Baz chan\" + bookmark + \"ges\" paragraph into two paragraphs:\n//    \"Baz changes\" (no bookmark) and a new paragraph \"Vi Du Luong changes \"\n//    that owns the _GoBack bookmark.\nif (bazPara) {\n  bazPara.insertText(\"Baz changes\", Word.InsertLocation.replace);\n  const newPara = bazPara.insertParagraph(\"Vi Du Luong changes \", Word.InsertLocation.after);\n  await context.sync();\n\n  // Re-find the text we just inserted so we can anchor the bookmark at the\n  // very end of it (after the trailing space, matching the source doc).\n  const results = newPara.search(\"Vi Du Luong changes \", { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    const match = results.items[0];\n    const endRange = match.getRange(Word.RangeLocation.end);\n    endRange.insertBookmark(\"_GoBack\");\n    await context.sync();\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the two paragraphs we need to touch by scanning their current text,\n# so the script is resilient to exact paragraph numbering.\n$gtIndex = 0\n$bazIndex = 0\nfor ($k = 1; $k -le $d.Paragraphs.Count; $k++) {\n    $t = $d.Paragraphs($k).Range.Text\n    if ($t -like \"*your stuff after this line*\") {\n        $gtIndex = $k\n    } elseif ($t -like \"*Baz chan*\") {\n        $bazIndex = $k\n    }\n}\n\n# 1) Normalize the \">>>  your stuff after this line >>>\" paragraph down to a\n#    single run (collapses the old proofing-error run split without changing\n#    the visible text). Delete-then-insert (rather than a plain Text=\n#    assignment) because the engine treats a same-looking Text= replace as a\n#    no-op when the final text matches the original.\nif ($gtIndex -gt 0) {\n    $p = $d.Paragraphs($gtIndex)\n    $full = $d.Range($p.Range.Start, $p.Range.End - 1)\n    $full.Delete()\n    $ins = $d.Range($p.Range.Start, $p.Range.Start)\n    $ins.InsertAfter(\">>>  your stuff after this line >>>\")\n}\n\n# 2) Split the \"Baz chan\" + bookmark + \"ges\" paragraph into two paragraphs:\n#    \"Baz changes\" (no bookmark) and a new paragraph \"Vi Du Luong changes \"\n#    that owns the _GoBack bookmark at its very end.\nif ($bazIndex -gt 0) {\n    $p = $d.Paragraphs($bazIndex)\n    $full = $d.Range($p.Range.Start, $p.Range.End - 1)\n    $full.Text = \"Baz changes\" + [char]13 + \"Vi Du Luong changes \"\n\n    $newPara = $d.Paragraphs($bazIndex + 1)\n\n    # Append a throw-away character first so the bookmark insertion point we\n    # compute next sits strictly *inside* the paragraph's text (not exactly\n    # on the paragraph-mark boundary) -- inserting a bookmark exactly at a\n    # freshly-created paragraph boundary relocates it to the wrong spot.\n    $endPos = $newPara.Range.End - 1\n    $tempIns = $d.Range($endPos, $endPos)\n    $tempIns.InsertAfter(\"#\")\n\n    $newPara2 = $d.Paragraphs($bazIndex + 1)\n    $bmPos = $newPara2.Range.End - 2\n    $bmRng = $d.Range($bmPos, $bmPos)\n    $d.Bookmarks.Add(\"_GoBack\", $bmRng)\n\n    $newPara3 = $d.Paragraphs($bazIndex + 1)\n    $tempPos = $newPara3.Range.End - 2\n    $tempRng = $d.Range($tempPos, $tempPos + 1)\n    $tempRng.Delete()\n}\n"}
